# Adjusted stack value score to emphasize value over total FPTS.
# New day's stacks were entered on the "Current" sheet (Tampa Bay Rays,
# Baltimore Orioles, Texas Rangers, Washington Nationals), which feeds the
# "RG table" rollup sheet, and the completed prior day's four stacks were
# appended to the "Season Log" history sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Current" - plug in the four new stacks for the day
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Current")

# Result markers for each of the four stacks (set first so the "RG table"
# rollup formulas that read them pick up the fresh values once the rest
# of each block's numbers are entered below).
$ws.Range("D8").Value = "Success"
$ws.Range("N8").Value = "Success"
$ws.Range("D17").Value = "Failure"

# --- Block 1: A1:D8 -> Tampa Bay Rays hitters (FD, DK) ---
$ws.Range("A1").Value = "Tampa Bay Rays hitters (FD, DK)"

$ws.Range("A3").Value = "Diaz"
$ws.Range("B3").Value = 3400
$ws.Range("C3").Value = 6.2

$ws.Range("A4").Value = "Garcia"
$ws.Range("B4").Value = 3500
$ws.Range("C4").Value = 22.2

$ws.Range("A5").Value = "Adames"
$ws.Range("B5").Value = 2500
$ws.Range("C5").Value = 32.2

$ws.Range("A6").Value = "Zunino"
$ws.Range("B6").Value = 2500
$ws.Range("C6").Value = 0

# --- Block 2: F1:I8 -> Baltimore Orioles hitters (FD, DK) ---
$ws.Range("F1").Value = "Baltimore Orioles hitters (FD, DK)"

$ws.Range("F3").Value = "Alberto"
$ws.Range("G3").Value = 2500
$ws.Range("H3").Value = 18.4

$ws.Range("F4").Value = "Nunez"
$ws.Range("G4").Value = 3600
$ws.Range("H4").Value = 6.2

$ws.Range("F5").Value = "Severino"
$ws.Range("G5").Value = 2800
$ws.Range("H5").Value = 59.6

$ws.Range("F6").Value = "Broxton"
$ws.Range("G6").Value = 2400
$ws.Range("H6").Value = 30.9

# --- Block 3: K1:N8 -> Texas Rangers hitters (FD, DK) ---
$ws.Range("K1").Value = "Texas Rangers hitters (FD, DK)"

$ws.Range("K3").Value = "Choo"
$ws.Range("L3").Value = 3900
$ws.Range("M3").Value = 31.4

$ws.Range("K4").Value = "DeShields"
$ws.Range("L4").Value = 2400
$ws.Range("M4").Value = 18.2

$ws.Range("K5").Value = "Cabrera"
$ws.Range("L5").Value = 3200
$ws.Range("M5").Value = 22.2

$ws.Range("K6").Value = "Odor"
$ws.Range("L6").Value = 2600
$ws.Range("M6").Value = 12.9

# --- Block 4: A10:D17 -> Washington Nationals hitters (FD, DK) ---
$ws.Range("A10").Value = "Washington Nationals hitters (FD, DK)"

$ws.Range("A12").Value = "Turner"
$ws.Range("B12").Value = 3800
$ws.Range("C12").Value = 15.4

$ws.Range("A13").Value = "Eaton"
$ws.Range("B13").Value = 3100
$ws.Range("C13").Value = 15.4

$ws.Range("A14").Value = "Adams"
$ws.Range("B14").Value = 3000
$ws.Range("C14").Value = 9.2

$ws.Range("A15").Value = "Dozier (INJ)"
$ws.Range("B15").Value = 2800
$ws.Range("C15").Value = 3

# Move the selection like the author did after entering the new stacks
$ws.Range("M9").Select()

# ---------------------------------------------------------------------
# Sheet "RG table" - selection moved while reviewing the new rollups
# ---------------------------------------------------------------------
$rg = $wb.Worksheets.Item("RG table")
$rg.Range("A1:A4").Select()

# ---------------------------------------------------------------------
# Sheet "Season Log" - append the results from the prior day's stacks
# ---------------------------------------------------------------------
$log = $wb.Worksheets.Item("Season Log")

# Finish off the existing Draftshot entry for 2019-05-30 (row 142)
$log.Range("C142").Value = "Cleveland Indians hitters (FD, DK)"
$log.Range("D142").Value = 3.34
$log.Range("E142").Value = "Failure"

# 2019-05-31 RG entries
$log.Range("A143").Value = 43616
$log.Range("B143").Value = "RG"
$log.Range("C143").Value = "Cleveland Indians hitters (FD, DK)"
$log.Range("D143").Value = 3.07
$log.Range("E143").Value = "Failure"

$log.Range("A144").Value = 43616
$log.Range("B144").Value = "RG"
$log.Range("C144").Value = "Boston Red Sox righties (FD, DK)"
$log.Range("D144").Value = 0.46
$log.Range("E144").Value = "Failure"

$log.Range("A145").Value = 43616
$log.Range("B145").Value = "RG"
$log.Range("C145").Value = "Los Angeles Angels hitters (FD, DK)"
$log.Range("D145").Value = 4.19
$log.Range("E145").Value = "Failure"

$log.Range("A146").Value = 43616
$log.Range("B146").Value = "RG"
$log.Range("C146").Value = "Colorado Rockies hitters (FD, DK)"
$log.Range("D146").Value = 9.08
$log.Range("E146").Value = "Success"

# 2019-06-01 Draftshot entry (no stack logged that day)
$log.Range("A147").Value = 43617
$log.Range("B147").Value = "Draftshot"

# 2019-06-02 RG entries
$log.Range("A148").Value = 43618
$log.Range("B148").Value = "RG"
$log.Range("C148").Value = "Colorado Rockies hitters (DK)"
$log.Range("D148").Value = 2.46
$log.Range("E148").Value = "Success"

$log.Range("A149").Value = 43618
$log.Range("B149").Value = "RG"
$log.Range("C149").Value = "Baltimore Orioles hitters (FD, DK)"
$log.Range("D149").Value = 0.27
$log.Range("E149").Value = "Failure"

# 2019-06-03 RG entry
$log.Range("A150").Value = 43619
$log.Range("B150").Value = "RG"
$log.Range("C150").Value = "Houston Astros hitters (FD, DK)"
$log.Range("D150").Value = 3.74
$log.Range("E150").Value = "Success"

# 2019-06-04 RG entries
$log.Range("A151").Value = 43620
$log.Range("B151").Value = "RG"
$log.Range("C151").Value = "Tampa Bay Rays hitters (FD, DK)"
$log.Range("D151").Value = 5.09
$log.Range("E151").Value = "Success"

$log.Range("A152").Value = 43620
$log.Range("B152").Value = "RG"
$log.Range("C152").Value = "Baltimore Orioles hitters (FD, DK)"
$log.Range("D152").Value = 10.19
$log.Range("E152").Value = "Success"

$log.Range("A153").Value = 43620
$log.Range("B153").Value = "RG"
$log.Range("C153").Value = "Texas Rangers hitters (FD, DK)"
$log.Range("D153").Value = 7
$log.Range("E153").Value = "Success"

$log.Range("A154").Value = 43620
$log.Range("B154").Value = "RG"
$log.Range("C154").Value = "Washington Nationals hitters (FD, DK)"
$log.Range("D154").Value = 3.39
$log.Range("E154").Value = "Failure"

# Reposition the frozen pane / selection like the author left it
$log.Application.ActiveWindow.ScrollRow = 140
$log.Range("G146").Select()

# Force a full recalculation so every dependent formula (RG table rollups,
# Season Log success/failure tallies) carries a fresh cached value.
$excel.Calculate()
